# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to match the freshly generated data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 0
$wsExpo.Range("F6").Value = 0
$wsExpo.Range("F7").Value = 0
$wsExpo.Range("F9").Value = 0
$wsExpo.Range("F10").Value = 0

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 0
$wsAll.Range("F3").Value = 109
$wsAll.Range("F5").Value = 0
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 406
$wsAll.Range("F10").Value = 467
